$wb = $excel.ActiveWorkbook

# --- Fill in the new Register sheet test-case rows (previously-empty cells) ---
$wsRegister = $wb.Worksheets.Item("Register")

# Row 2: "Register 1" test case
$wsRegister.Range("B2").Value = "Weryfikacja formatki"
$wsRegister.Range("C2").Value = "Łukaś"
$wsRegister.Range("D2").Value = "Done"
$wsRegister.Range("E2").Value = "brak"

# Row 3: "Register 2" test case
$wsRegister.Range("B3").Value = "Rejestracja użytkownika - prawidłowe dane"
$wsRegister.Range("C3").Value = "Łukaś"
$wsRegister.Range("D3").Value = "Done"
$wsRegister.Range("E3").Value = "brak"

# --- Update selections / active sheet so Register becomes the visible tab ---
$wsLogin = $wb.Worksheets.Item("Login")
$wsLogin.Activate() | Out-Null
$wsLogin.Range("A1:E4").Select() | Out-Null

$wsRegister.Activate() | Out-Null
$wsRegister.Range("E3").Select() | Out-Null

# --- Best-effort cosmetic cleanup (locale rename of the built-in style) ---
try {
    $wb.Styles.Item("Normalny").Name = "Normal"
} catch {
}
